$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23/24: Coin/Link content swapped (Polygon <-> Dai)
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"

# Price column (D) updates - force text storage (matches source data, e.g. "67.414.89")
$priceUpdates = @{
    "D2"  = "67.414.89";
    "D3"  = "3.527.85";
    "D5"  = "596.83";
    "D6"  = "174.09";
    "D7"  = "1.00";
    "D8"  = "0.593";
    "D10" = "7.33";
    "D11" = "0.438";
    "D12" = "4.140.61";
    "D14" = "28.82";
    "D16" = "67.286.33";
    "D17" = "3.544.43";
    "D18" = "6.38";
    "D19" = "14.24";
    "D20" = "397.21";
    "D21" = "8.02";
    "D22" = "73.78";
    "D23" = "1.00";
    "D24" = "0.541";
    "D25" = "0.0000125";
    "D26" = "10.30";
    "D28" = "0.998";
    "D29" = "6.35";
    "D31" = "2.08";
    "D32" = "24.12";
    "D33" = "7.45";
    "D35" = "163.47";
    "D36" = "0.903";
    "D39" = "6.83";
    "D40" = "0.0751";
    "D41" = "26.72";
    "D42" = "27.40";
    "D43" = "2.65";
    "D44" = "2.816.21";
    "D45" = "42.97";
    "D47" = "344.15";
    "D49" = "33.79";
    "D50" = "0.861";
    "D51" = "6.55";
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}

# Volume(1h) column (E) updates
$volumeUpdates = @{
    "E2"  = "  +0.98%  ";
    "E3"  = "  +0.83%  ";
    "E4"  = "  +0.02%  ";
    "E5"  = "  +0.77%  ";
    "E6"  = "  +2.12%  ";
    "E7"  = "  +0.01%  ";
    "E8"  = "  +1.45%  ";
    "E9"  = "  +8.07%  ";
    "E10" = "  +0.80%  ";
    "E11" = "  +0.32%  ";
    "E12" = "  +0.99%  ";
    "E14" = "  +2.66%  ";
    "E15" = "  +2.67%  ";
    "E16" = "  +0.88%  ";
    "E17" = "  +1.37%  ";
    "E18" = "  +1.43%  ";
    "E19" = "  +1.19%  ";
    "E20" = "  +2.09%  ";
    "E21" = "  +0.47%  ";
    "E22" = "  +1.05%  ";
    "E23" = "  +0.29%  ";
    "E24" = "  +2.49%  ";
    "E25" = "  +0.14%  ";
    "E26" = "  +1.02%  ";
    "E27" = "  +0.16%  ";
    "E28" = "  -0.25%  ";
    "E29" = "  -0.43%  ";
    "E30" = "  -0.10%  ";
    "E31" = "  +1.31%  ";
    "E32" = "  +2.69%  ";
    "E33" = "  +0.62%  ";
    "E34" = "  +5.21%  ";
    "E35" = "  +1.55%  ";
    "E36" = "  +0.01%  ";
    "E37" = "  -0.33%  ";
    "E38" = "  +3.32%  ";
    "E39" = "  +1.63%  ";
    "E40" = "  +0.90%  ";
    "E41" = "  +1.24%  ";
    "E42" = "  +0.80%  ";
    "E43" = "  +2.97%  ";
    "E44" = "  +0.42%  ";
    "E45" = "  -1.32%  ";
    "E46" = "  -0.54%  ";
    "E47" = "  -3.62%  ";
    "E48" = "  +1.16%  ";
    "E49" = "  +2.39%  ";
    "E50" = "  +1.06%  ";
    "E51" = "  +1.29%  ";
}

foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
